$d = $word.ActiveDocument

function Find-ParagraphIndex($doc, [string]$needle) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        if ($doc.Paragraphs($i).Range.Text.Contains($needle)) {
            return $i
        }
    }
    return -1
}

# ---------------------------------------------------------------------------
# 1) "CORRIGIR ERRO DE LÓGICA DA MOVIMENTAÇÃO DO CARRO DA POLICIA NO SENTIDO
#    DIREITO" -> split into "CORRIGIR " + "SISTEMA DE ESQUIVA DO PLAYER PARA
#    QUE CONTAGEM DE PONTOS NÃO SEJA VÁLIDA AO OCORRER COLISÃO"
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "CORRIGIR ERRO DE LÓGICA DA MOVIMENTAÇÃO DO CARRO DA POLICIA NO SENTIDO DIREITO",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "CORRIGIR SISTEMA DE ESQUIVA DO PLAYER PARA QUE CONTAGEM DE PONTOS NÃO SEJA VÁLIDA AO OCORRER COLISÃO",
    2) | Out-Null

$idx1 = Find-ParagraphIndex $d "CORRIGIR SISTEMA DE ESQUIVA DO PLAYER"
$p1 = $d.Paragraphs($idx1)
$p1Start = $p1.Range.Start
$p1End = $p1.Range.End - 1

# Split "CORRIGIR " off into its own run (toggling a format on/off forces the
# engine to materialize a run boundary without leaving residual formatting).
$p1RunA = $d.Range($p1Start, $p1Start + 9)
$p1RunA.Bold = 1
$p1RunA.Bold = 0

$p1RunB = $d.Range($p1Start + 9, $p1End)
$p1RunB.Bold = 1
$p1RunB.Bold = 0

# ---------------------------------------------------------------------------
# 2) "- CARRO NÃO TÁ BLOQUEANDO PLAYER" -> split into "- " + "ESTÁ SENDO
#    CONTADO PONTOS DE ESQUIVA MESMO QUANDO PLAYER COLIDE" (leading space run
#    stays as-is)
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "- CARRO NÃO TÁ BLOQUEANDO PLAYER",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "- ESTÁ SENDO CONTADO PONTOS DE ESQUIVA MESMO QUANDO PLAYER COLIDE",
    2) | Out-Null

$idx2 = Find-ParagraphIndex $d "ESTÁ SENDO CONTADO PONTOS DE ESQUIVA"
$p2 = $d.Paragraphs($idx2)
$p2Start = $p2.Range.Start
$p2End = $p2.Range.End - 1

$p2RunA = $d.Range($p2Start + 1, $p2Start + 3)
$p2RunA.Bold = 1
$p2RunA.Bold = 0

$p2RunB = $d.Range($p2Start + 3, $p2End)
$p2RunB.Bold = 1
$p2RunB.Bold = 0
